$wb = $excel.ActiveWorkbook

# --- 1. Insert a new "State" column into hotel_info, right after "Hotel_Name" ---
$hotelSheet = $wb.Worksheets.Item("hotel_info")

# Column C currently holds "City" (B=Hotel_Name). Insert a new column there.
$hotelSheet.Columns.Item(3).Insert()

$hotelSheet.Cells.Item(1, 3).Value = "State"
$hotelSheet.Cells.Item(2, 3).Value = "Louisiana"

# Insert() duplicated the shifted row into the new empty Local_Rank cell
# (now column I) with a bogus value - clear it back to blank.
$hotelSheet.Cells.Item(2, 9).Value = ""

# --- 2. Reorder the sheets: review_info first, hotel_info second ---
$reviewSheet = $wb.Worksheets.Item("review_info")
$reviewSheet.Move($hotelSheet)
